$d = $word.ActiveDocument

# Mapping of old font sizes (points) to new font sizes (points), per the
# commit: name 16->18, contact/body 9->10, section headers 12->13,
# job titles 11->12, overview/section paragraphs 10->11.
$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    11 = 12
    10 = 11
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    # Exclude the trailing paragraph mark so only the run(s) with actual
    # text get the new size written to their <w:rPr>, instead of also
    # stamping a <w:rPr> onto the paragraph mark (<w:pPr><w:rPr>).
    $r = $d.Range($full.Start, $full.End - 1)
    if ($r.Start -ge $r.End) {
        continue
    }
    $current = $r.Font.Size
    if ($sizeMap.ContainsKey($current)) {
        $r.Font.Size = $sizeMap[$current]
    }
}
